$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab19")

$ws.Range("R3").Value = 45.025733623187897
$ws.Range("R4").Value = 41.891955685351803
$ws.Range("R5").Value = 45.040888826404696
$ws.Range("R6").Value = 44.345664757404002
$ws.Range("R7").Value = 10.3727474306876
$ws.Range("R8").Value = 40.5940552908055
$ws.Range("R9").Value = 34.331492876482699
$ws.Range("R10").Value = 31.193814062404901
$ws.Range("R11").Value = 52.951540056276201
$ws.Range("R12").Value = 18.279911392911501
$ws.Range("R13").Value = 34.5944064075684
$ws.Range("R14").Value = 8.1227926531516701
$ws.Range("R15").Value = 16.4416262007732
$ws.Range("R16").Value = 13.688864979750299
$ws.Range("R17").Value = 34.529957898005001
$ws.Range("R18").Value = 72.891877155626503
$ws.Range("R19").Value = 39.553458604219401
$ws.Range("R20").Value = 44.622731725259897
$ws.Range("R21").Value = 39.862196110500797
$ws.Range("R22").Value = 15.309847129227901
$ws.Range("R23").Value = 33.056170556374198
$ws.Range("R24").Value = 11.516798569849399
$ws.Range("R25").Value = 149.939648857761
$ws.Range("R26").Value = 38.753396744285503
$ws.Range("R27").Value = 8.5738254450976097
$ws.Range("R28").Value = 10.7162886540966
$ws.Range("R29").Value = 23.032116454966602
$ws.Range("R30").Value = 27.8335084103309
$ws.Range("R31").Value = 19.064319772189801
$ws.Range("R32").Value = 120.178646783612
$ws.Range("R33").Value = 17.178050909260499
$ws.Range("R34").Value = 71.500734394283
$ws.Range("R35").Value = 17.0172856140736
$ws.Range("R36").Value = 14.1751950000597
$ws.Range("R37").Value = 14.385308686659201
$ws.Range("R38").Value = 14.6283437020997
$ws.Range("R39").Value = 25.605455733878198
$ws.Range("R40").Value = 10.5531691004687
$ws.Range("R42").Value = 33.2069791354443
$ws.Range("R43").Value = 32.985864623713397
$ws.Range("R44").Value = 42.439990494046803
$ws.Range("R45").Value = 18.783189127890399
$ws.Range("R46").Value = 23.482995586326499
$ws.Range("R47").Value = 29.4574859455408
$ws.Range("R48").Value = 28.9751184139811
$ws.Range("R49").Value = 22.632864580266698
$ws.Range("R50").Value = 9.1788319187919996
$ws.Range("R51").Value = 30.1948244523951
$ws.Range("R52").Value = 63.545103442808397
$ws.Range("R53").Value = 19.4231713892963
$ws.Range("R54").Value = 32.0125450422886
$ws.Range("R55").Value = 27.3971836887204
$ws.Range("R56").Value = 15.496881969881301
$ws.Range("R57").Value = 11.499661296247799
$ws.Range("R58").Value = 23.373588354331201
$ws.Range("R59").Value = 26.864112802491
$ws.Range("R60").Value = 22.5931424248264
$ws.Range("R61").Value = 17.668744907026699
$ws.Range("R62").Value = 21.428809165303299
$ws.Range("R63").Value = 29.074255902749801
$ws.Range("R64").Value = 27.112143836084599
$ws.Range("R65").Value = 24.0965283884358
$ws.Range("R66").Value = 28.688994384145701
$ws.Range("R67").Value = 15.827399720297199
$ws.Range("R68").Value = 17.5819635585536
$ws.Range("R69").Value = 13.418753152565801
$ws.Range("R70").Value = 36.685339214373897
$ws.Range("R71").Value = 17.668744907026699
$ws.Range("R72").Value = 13.3876320778232
$ws.Range("R73").Value = 32.443122353982901
$ws.Range("R74").Value = 30.495426343310001
$ws.Range("R75").Value = 43.522693584562703
$ws.Range("R76").Value = 52.997166297673402
$ws.Range("R77").Value = 21.190808658086699
$ws.Range("R78").Value = 50.405171963548
$ws.Range("R79").Value = 30.026167846821199
$ws.Range("R80").Value = 20.986150259167299
$ws.Range("R81").Value = 37.656844581745602
$ws.Range("R82").Value = 21.601179840296101
$ws.Range("R83").Value = 28.295846675813898
$ws.Range("R84").Value = 21.1779638950382
$ws.Range("R85").Value = 12.080169759101199
$ws.Range("R86").Value = 19.234074729172399
$ws.Range("R87").Value = 25.213955396743401
$ws.Range("R88").Value = 32.2537412319464
$ws.Range("R89").Value = 26.486034597540002
$ws.Range("R90").Value = 32.383493587139697
$ws.Range("R91").Value = 25.363513289198099
$ws.Range("R92").Value = 15.369241775483101
$ws.Range("R93").Value = 32.089359681950697
$ws.Range("R94").Value = 117.52915593533599
$ws.Range("R95").Value = 20.347630778781198
$ws.Range("R96").Value = 31.3023002992424
$ws.Range("R97").Value = 19.850907550630101
$ws.Range("R98").Value = 21.334067301334301
